# Commit: "Se elimina todo de flyio para cambiar a heroku"
# The only functional content change touching this workbook is a text fix in
# the "Importar Propiedades" sheet header row: the accented word
# "localización" is replaced with the unaccented "localizacion".
#
# (The rest of the original diff -- shifted cellXfs indices, a cloned/near-
# duplicate style inserted at index 1, font "charset" attributes, every xf's
# applyProtection flag flipping to true, and a hundredth-of-a-character
# rounding tweak to a column width on the second sheet -- are exactly the
# kind of inconsequential re-serialization noise a spreadsheet engine leaves
# behind when it resaves a workbook that was originally produced by a
# different writer; they don't change any cell's visible formatting. They
# are reproduced here only where doing so is meaningful/controllable --
# e.g. restoring the active-cell selection -- and otherwise left to the
# host application's own normalization on save.)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Fix the accented header text.
$ws1.Range("B1").Value = "localizacion"

# Restore the active selection recorded in the sheet view.
$ws1.Range("B5").Select()
